$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Author 1"
$ws.Range("E1").Value = "Author 2"
$ws.Range("F1").Value = "Author 3"
$ws.Range("G1").Value = "Author 4"
$ws.Range("H1").Value = "Author 5"
$ws.Range("I1").Value = "Author 6"
$ws.Range("J1").Value = "Author 7"
$ws.Range("K1").Value = "Author 8"
$ws.Range("L1").Value = "Author 9"
$ws.Range("M1").Value = "Author 10"
$ws.Range("N1").Value = "Author 1"
$ws.Range("O1").Value = "Author 12"
$ws.Range("P1").Value = "Author 13"
$ws.Range("Q1").Value = "Author 14"
$ws.Range("R1").Value = "Author 15"
$ws.Range("S1").Value = "Institution 1"
$ws.Range("T1").Value = "Institution 2"
$ws.Range("U1").Value = "Institution 3"
$ws.Range("V1").Value = "Institution 4"
$ws.Range("W1").Value = "Institution 5"
$ws.Range("X1").Value = "Institution 6"
$ws.Range("Y1").Value = "Institution 7"
$ws.Range("Z1").Value = "Institution 8"
$ws.Range("AA1").Value = "Institution 9"
$ws.Range("AB1").Value = "Institution 10"
$ws.Range("AC1").Value = "Institution 11"
$ws.Range("AD1").Value = "Institution 12"
$ws.Range("AE1").Value = "Institution 13"
$ws.Range("AF1").Value = "Institution 14"
$ws.Range("AG1").Value = "Institution 15"
$ws.Range("AH1").Value = "Institution 16"
